# Fehler.xlsx ("Nicht Programm") – bug tracker update
#  - fixed the wording of the existing "Auto Btn" KI bug entry (row 3)
#  - added the new "KivsKi" ChoiceBox event bug (row 15)
#  - added the new PutShips-Multiplayer-hang bug (row 16)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 3: reword the existing bug description (button bug fix) ---
$ws.Range("A3").Value = "Auto Btn im KI Multiplayer Spiel führt zu Fehlern in der Kommunikation tut noch nicht"
$ws.Range("A3").WrapText = $true

# --- Row 15: new bug entry (KivsKi ChoiceBox / framework for kistrong) ---
$ws.Range("A15").Value = "Falsches ChoiceBox Event in KivsKi, Geschwindigkeit updated beim Klicken auf die Box und nicht auf das Menü"
$ws.Range("A15").WrapText = $true
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Controller_GameScreen"
$ws.Range("D15").Value = "einfach"

# --- Row 16: new bug entry (PutShips hangs in Multiplayer) ---
$ws.Range("A16").Value = "PutShips hängt sich im Multiplayer auf wenn man vor dem Host auf Start drückt"
$ws.Range("A16").WrapText = $true
